# Edit the CasesTab query text in cell B2: remove the trailing "Cohort" column
# from the RETURN clause (fixed extra browser opening in each script; added implicit wait).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCasesTabQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC01'] and demo.breed in ['Irish Terrier','Poodle','Wire Fox Terrier']and diag.disease_term in ['Bladder Cancer'] and diag.primary_disease_site in ['Bladder, Urethra']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$ws.Range("B2").Value2 = $newCasesTabQuery

# Row heights shrink slightly now that the query text is one line shorter.
$ws.Rows.Item(2).RowHeight = 300
$ws.Rows.Item(3).RowHeight = 285
$ws.Rows.Item(4).RowHeight = 270

# Match the saved selection state.
$ws.Range("B3").Select()
